$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.270.69'
$ws.Range("E2").Value = '  +2.08%  '
$ws.Range("D3").Value = '2.375.90'
$ws.Range("E3").Value = '  +0.00%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("E5").Value = '  +8.25%  '
$ws.Range("D6").Value = "'244.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.37%  '
$ws.Range("D7").Value = "'76.67"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +5.80%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = "'0.598"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +24.28%  '
$ws.Range("D10").Value = "'0.104"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.25%  '
$ws.Range("E11").Value = '  +1.79%  '
$ws.Range("D12").Value = "'32.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +19.20%  '
$ws.Range("D13").Value = "'7.51"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +19.14%  '
$ws.Range("E14").Value = '  +2.61%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = "'17.23"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +7.00%  '
$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C16").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D16").Value = '2.731.51'
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("D17").Value = "'0.927"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +7.73%  '
$ws.Range("D18").Value = '2.378.54'
$ws.Range("E18").Value = '  +0.23%  '
$ws.Range("D19").Value = '44.242.85'
$ws.Range("E19").Value = '  +1.98%  '
$ws.Range("E20").Value = '  +4.37%  '
$ws.Range("D21").Value = "'6.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.85%  '
$ws.Range("D22").Value = "'78.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.57%  '
$ws.Range("D23").Value = "'258.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.30%  '
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("D25").Value = "'2.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.82%  '
$ws.Range("D26").Value = "'3.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.65%  '
$ws.Range("D27").Value = "'10.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +8.01%  '
$ws.Range("D28").Value = "'1.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +14.61%  '
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").Value = "'23.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.07%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = "'2.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.13%  '
$ws.Range("D31").Value = "'175.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.76%  '
$ws.Range("E32").Value = '  +1.12%  '
$ws.Range("E33").Value = '  +6.63%  '
$ws.Range("D34").Value = "'5.41"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +8.03%  '
$ws.Range("D35").Value = "'0.0765"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +10.51%  '
$ws.Range("E36").Value = '  +5.28%  '
$ws.Range("E37").Value = '  +5.76%  '
$ws.Range("D38").Value = "'2.50"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.43%  '
$ws.Range("E39").Value = '  -0.93%  '
$ws.Range("D40").Value = "'0.0278"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.36%  '
$ws.Range("B41").Value = 'InjectiveProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D41").Value = "'19.25"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.03%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = "'9.15"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.46%  '
$ws.Range("E43").Value = '  +0.05%  '
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").Value = "'0.195"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +16.47%  '
$ws.Range("B45").Value = 'ARBITRUM'
$ws.Range("C45").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D45").Value = "'1.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.39%  '
$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D46").Value = "'1.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.74%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = "'0.101"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.65%  '
$ws.Range("D48").Value = "'2.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +12.39%  '
$ws.Range("D49").Value = "'102.51"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.94%  '
$ws.Range("D50").Value = "'4.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.22%  '
$ws.Range("D51").Value = '1.481.94'
$ws.Range("E51").Value = '  +2.14%  '
